# Truncate the decimal values in column L ("txp/alvl/adist") of the
# distance_quest_data sheet down to whole numbers (drop the fractional part),
# for every data row of the table (rows 2 through 130).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 130
$col = 12  # column L

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $v = $cell.Value2
    if ($v -ne $null) {
        $truncated = [int]$v
        if ($truncated -ne $v) {
            $cell.Value2 = $truncated
        }
    }
}
